$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the row that held the value 43 (it is being folded into the
#        "43+1" label placed in B30); this shifts rows 32-39 up by one,
#        leaving B31:B38 = 44..51 and removing the old row 39.
$ws.Rows.Item(31).Delete()

# --- 2. Give the header row (A1:B1) a new fill color (Accent4 theme color).
#        First "prime" the cells with an existing solid pattern fill (copied
#        from B26, which already uses a solid fill) so that only the final
#        color needs a new style entry.
$ws.Range("B26").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("A1:B1").Interior.ThemeColor = 8

# --- 3. B28:B30 adopt the same highlight style already used by B27.
$ws.Range("B27").Copy()
$ws.Range("B28:B30").PasteSpecial(-4122)

# --- 4. Update the values: B28 and B30 become text labels ("40+3", "43+1"),
#        B29 keeps its original numeric value (41).
$ws.Range("B28").Value = "40+3"
$ws.Range("B30").Value = "43+1"

# --- 5. Update the selection / scrolled view to match the new state.
$ws.Range("F21").Select()
